$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text corrections to the use-case wording ---

# Step 3 (D9): "cronómetro" -> "temporizador"
$ws.Range("D9").Value = "3. Inicia temporizador para o passo"

# Step 4 (D11): "cronómetro" -> "temporizador" and clarify that the value gets recorded
$ws.Range("D11").Value = "4. Termina temporizador para o passo e regista o seu valor"

# Alternativa 3 steps were renumbered: old 3.1/3.2 collapse into a single 4.1 step,
# and the old 3.2 row is cleared out
$ws.Range("D22").Value = "4.1 Regressa a 1"
$ws.Range("D23").Value = ""

# Alternativa 5 (D28): add the missing "para" before "fim da contagem"
$ws.Range("D28").Value = "(1/2/3/4).1 Alerta para fim da contagem e apresenta sugestão de ação"

# --- Interface / formatting fixes ---

# D15 was missing word-wrap even though its neighbours (D19/D23/D26) wrap long text
$ws.Range("D15").WrapText = $true

# D11 was missing the left/right border used by the rest of column D (it only had
# the default no-border style) - match it to D15/D19/D23/D26's border
$ws.Range("D11").Borders.Item(7).LineStyle = 1
$ws.Range("D11").Borders.Item(7).Weight = -4138
$ws.Range("D11").Borders.Item(10).LineStyle = 1
$ws.Range("D11").Borders.Item(10).Weight = -4138

# Reflect the last-active cell after these edits (matches the saved selection)
$ws.Range("D23").Select()
